# Auto commit at 2025-09-20 8:08:09.05
#
# Updates the "Metrics" sheet's monthly/annual/total figures and moves the
# active tab / selection from "Metrics" to "today" (the "today" sheet's
# figures are plain formulas referencing Metrics!B2:B13, so they (and the
# TODAY()-1 date cell) recalculate automatically).

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# --- Updated metric values (Metrics!B2:B13) ---------------------------------
$wsMetrics.Range("B2").Value  = 298308.58999999997   # month-charge
$wsMetrics.Range("B3").Value  = 241255.47000000003   # month-electricity
$wsMetrics.Range("B4").Value  = 94248.960000000006   # month-service
$wsMetrics.Range("B5").Value  = 11694                # month-orders
$wsMetrics.Range("B6").Value  = 4217559.47           # year-charge
$wsMetrics.Range("B7").Value  = 3568782.9499999993   # year-electricity
$wsMetrics.Range("B8").Value  = 1223614.6400000001   # year-service
$wsMetrics.Range("B9").Value  = 162854               # year-orders
$wsMetrics.Range("B10").Value = 32682883.270999826   # total-charge
$wsMetrics.Range("B11").Value = 19598653.020000003   # total-electricity
$wsMetrics.Range("B12").Value = 11505323.530000001   # total-service
$wsMetrics.Range("B13").Value = 1260481              # total-orders

# --- Selection on the Metrics sheet moves from H12 to G11 -------------------
$wsMetrics.Range("G11").Select() | Out-Null

# --- "today" becomes the active tab, selection moves from E7 to F7 ----------
$wsToday.Range("F7").Select() | Out-Null
$wsToday.Activate() | Out-Null
